$d = $word.ActiveDocument

# The meeting-minutes table is the first (only) table in the document.
$table = $d.Tables.Item(1)

# The last row is the empty one awaiting the newest entry.
$lastRow = $table.Rows.Item($table.Rows.Count)

# First cell (Date column) -> "28/12"
$lastRow.Cells.Item(1).Range.Text = "28/12"

# Second cell (Description column) -> pairing note
$lastRow.Cells.Item(2).Range.Text = "Pair programming while building [PAG3]"
